$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.895.93"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.640.25"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.70"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2605"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06468"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07814"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.660.11"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.265"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "1.866.18"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5660"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "0.0₅7714"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.53"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "25.905.16"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.90"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.992"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.228"
$ws.Range("E23").Value = "  +4.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.767"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.20"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1233"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.869"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.245"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05031"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.315"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.258"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.579"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.385"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9088"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.579"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5536"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "1.131.62"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9953"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.00"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.500"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8020"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.82"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4235"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.749"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05046"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +0.10%  "
